{"js": "// Update codes for recalibration\n// Replace each reported estimate (%) and its 95% CI text with the\n// recalibrated value throughout the document (table cells).\n\nconst replacements = [\n    [\"1.92%\", \"2.08%\"],\n    [\" (1.62% to 2.21%)\", \" (1.8% to 2.37%)\"],\n    [\"4.79%\", \"4.33%\"],\n    [\" (4.32% to 5.24%)\", \" (3.92% to 4.74%)\"],\n    [\"3.45%\", \"3.79%\"],\n    [\" (3.05% to 3.86%)\", \" (3.4% to 4.18%)\"],\n    [\"9.48%\", \"8.88%\"],\n    [\" (8.83% to 10.12%)\", \" (8.3% to 9.47%)\"],\n    [\"4.47%\", \"5.02%\"],\n    [\" (4.01% to 4.93%)\", \" (4.57% to 5.48%)\"],\n    [\"14.38%\", \"13.73%\"],\n    [\" (13.59% to 15.16%)\", \" (13.01% to 14.45%)\"],\n    [\"5.33%\", \"5.93%\"],\n    [\" (4.82% to 5.84%)\", \" (5.44% to 6.43%)\"],\n    [\"19.96%\", \"19.09%\"],\n    [\" (19.04% to 20.87%)\", \" (18.24% to 19.92%)\"],\n    [\"6.06%\", \"6.74%\"],\n    [\" (5.51% to 6.61%)\", \" (6.2% to 7.27%)\"],\n    [\"25.64%\", \"24.92%\"],\n    [\" (24.59% to 26.67%)\", \" (23.96% to 25.87%)\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n    const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    await context.sync();\n\n    for (let i = 0; i < results.items.length; i++) {\n        results.items[i].insertText(newText, Word.InsertLocation.replace);\n    }\n    await context.sync();\n}\n", "ps1": "# Update codes for recalibration\n# Replace each reported estimate (%) and its 95% CI text with the\n# recalibrated value throughout the document (table cells).\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"1.92%\", \"2.08%\"),\n    @(\" (1.62% to 2.21%)\", \" (1.8% to 2.37%)\"),\n    @(\"4.79%\", \"4.33%\"),\n    @(\" (4.32% to 5.24%)\", \" (3.92% to 4.74%)\"),\n    @(\"3.45%\", \"3.79%\"),\n    @(\" (3.05% to 3.86%)\", \" (3.4% to 4.18%)\"),\n    @(\"9.48%\", \"8.88%\"),\n    @(\" (8.83% to 10.12%)\", \" (8.3% to 9.47%)\"),\n    @(\"4.47%\", \"5.02%\"),\n    @(\" (4.01% to 4.93%)\", \" (4.57% to 5.48%)\"),\n    @(\"14.38%\", \"13.73%\"),\n    @(\" (13.59% to 15.16%)\", \" (13.01% to 14.45%)\"),\n    @(\"5.33%\", \"5.93%\"),\n    @(\" (4.82% to 5.84%)\", \" (5.44% to 6.43%)\"),\n    @(\"19.96%\", \"19.09%\"),\n    @(\" (19.04% to 20.87%)\", \" (18.24% to 19.92%)\"),\n    @(\"6.06%\", \"6.74%\"),\n    @(\" (5.51% to 6.61%)\", \" (6.2% to 7.27%)\"),\n    @(\"25.64%\", \"24.92%\"),\n    @(\" (24.59% to 26.67%)\", \" (23.96% to 25.87%)\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
